$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("23:25").Insert()
$ws.Range("A22:J25").FillDown()
Write-Host "done"
